# 4/12 Added hover effects to buttons and general UI Polish
#
# D1/E1/F1 on the active sheet get updated to the text value "242"
# (previously "204", "204" and "51" respectively). Because "242" looks
# like a number, Excel would auto-convert a plain .Value assignment to a
# numeric cell, so we briefly mark the range as Text, type the values in,
# then clear the formatting back off so the cells keep their original
# (default) style while still storing "242" as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("D1:F1")
$target.NumberFormat = "@"

$ws.Range("D1").Value = "242"
$ws.Range("E1").Value = "242"
$ws.Range("F1").Value = "242"

$target.ClearFormats()
